# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-13 05:21:19
#
# The "Recorded By" column (G) lists the accounts that touched a session's
# attendance record as a comma-separated string. Upstream regenerated the
# report and the ordering of that list shifted: the first entry in each
# cell's list was moved to the end (a left-rotation by one element).
# Reproduce that same rotation here, cell by cell, for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Text

    if ([string]::IsNullOrEmpty($orig)) {
        continue
    }

    $rawParts = $orig.Split(",")
    if ($rawParts.Count -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $newVal = [string]::Join(", ", $rotated)

    if ($newVal -ne $orig) {
        $cell.Value = $newVal
    }
}
